$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76, shifting existing rows 76:82 down to 77:83
$ws.Rows("76:76").Insert()

# Populate the newly inserted row 76 with the new weekly data point
$ws.Cells.Item(76, 1).Value = 11
$ws.Cells.Item(76, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(76, 3).Value = "Bíobío"
$ws.Cells.Item(76, 4).Value = 45132
$ws.Cells.Item(76, 5).Value = 8
$ws.Cells.Item(76, 6).Value = 100112043
$ws.Cells.Item(76, 7).Value = "Pepino dulce"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 100
$ws.Cells.Item(76, 11).Value = 16000
$ws.Cells.Item(76, 12).Value = 17000
$ws.Cells.Item(76, 13).Value = 16500
$ws.Cells.Item(76, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(76, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(76, 16).Value = 917
$ws.Cells.Item(76, 17).Value = 18
$ws.Cells.Item(76, 18).Value = "Hortaliza"
